# Update the "Minimum historic data required" value for the "Mean" algorithm
# row from "1 year" to "26 weeks".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D5").Value = "26 weeks"
